$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.401.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.236.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.87%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("E13").Value = "  -0.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.853"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.260.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.211.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("E18").Value = "  +11.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +38.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.02%  "

$ws.Range("E24").Value = "  +3.44%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0808"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.118"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.97%  "

$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0309"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.82%  "

$ws.Range("E39").Value = "  -0.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.201"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.16%  "

$ws.Range("E43").Value = "  +1.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.32%  "

$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.06%  "

$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.26%  "

